# Updating to reflect new gtest
#
# The "BST Search" rubric row is retired (merged into a single
# "BST SearchFoundAndNotFound" test), so:
#   1. Delete the "BST Search" row (row 7) from the Classroom sheet -
#      this naturally shifts every row below it up by one and keeps
#      the Total Points SUM formula's range in sync.
#   2. The row that used to be "BST SearchNotFound" (now row 7) is
#      renamed to "BST SearchFoundAndNotFound" and its point value
#      changes from 10 to 40.
#   3. Update the active selection to A13 to match the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Classroom")

# 1. Remove the obsolete "BST Search" row entirely.
$ws.Rows(7).Delete() | Out-Null

# 2. Rename/repoint what is now row 7 and update its points.
$ws.Range("A7").Value = "BST SearchFoundAndNotFound"
$ws.Range("E7").Value = 40

# 3. Match the saved selection/active cell.
$ws.Range("A13").Select() | Out-Null
